# refcode 2 updating when reupload stuff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: make room for a "Reference Code 2" column right after "Reference Code" ---
# Insert a blank column at G (this is what shifts the existing column-width
# customisations, originally on columns J/O, out to K/P - a plain Cut+Insert of
# the old Reference Code 2 column would drag its own width along instead).
$ws.Columns("G:G").Insert()
# The insert pushed the old "Reference Code" column out to Q; nothing belongs
# there any more since we rewrite every header/value below.
$ws.Columns("Q:Q").ClearContents()

# --- Step 2: rewrite the header row with "Reference Code 2" in its new slot ---
$ws.Range("A1").Value = "Receipt ID"
$ws.Range("B1").Value = "Date"
$ws.Range("C1").Value = "Amount"
$ws.Range("D1").Value = "Recurring Total Months"
$ws.Range("E1").Value = "Recurrence Number"
$ws.Range("F1").Value = "Donor First Name"
$ws.Range("G1").Value = "Reference Code 2"
$ws.Range("H1").Value = "Donor Last Name"
$ws.Range("I1").Value = "Donor Addr1"
$ws.Range("J1").Value = "Donor City"
$ws.Range("K1").Value = "Donor State"
$ws.Range("L1").Value = "Donor ZIP"
$ws.Range("M1").Value = "Donor Country"
$ws.Range("N1").Value = "Donor Email"
$ws.Range("O1").Value = "Donor Phone"
$ws.Range("P1").Value = "Reference Code"

# --- Step 3: overwrite rows 2-4 with the new donor data (reupload) ---

# Row 2 - Reba Wooden
$ws.Range("A2").Value = "AB191624287"
$ws.Range("B2").Value2 = 44338.17931712963
$ws.Range("C2").Value = 50
$ws.Range("D2").Value = "unlimited"
$ws.Range("E2").Value = 4
$ws.Range("F2").Value = "Reba"
$ws.Range("G2").Value = "TEST 9"
$ws.Range("H2").Value = "Wooden"
$ws.Range("I2").Value = "113 Severn Drive"
$ws.Range("J2").Value = "Greenwood"
$ws.Range("K2").Value = "IN"
$ws.Range("L2").Value = 46142
$ws.Range("M2").Value = "United States"
$ws.Range("N2").Value = "rboydw@gmsil.com"
$ws.Range("O2").Value = 3177975892
$ws.Range("P2").Value = "2.22.21.EOM1."

# Row 3 - Linda Braun
$ws.Range("A3").Value = "AB191627990"
$ws.Range("B3").Value2 = 44338.179351851853
$ws.Range("C3").Value = 25
$ws.Range("D3").Value = "unlimited"
$ws.Range("E3").Value = 4
$ws.Range("F3").Value = "Linda"
$ws.Range("G3").ClearContents()
$ws.Range("H3").Value = "Braun"
$ws.Range("I3").Value = "11 Parker Road"
$ws.Range("J3").Value = "Arlington"
$ws.Range("K3").Value = "MA"
$ws.Range("L3").Value = 2474
$ws.Range("M3").Value = "United States"
$ws.Range("N3").Value = "lbraun@verizon.net"
$ws.Range("O3").Value = 7816462999
$ws.Range("P3").Value = "2.22.21.EOM1."

# Row 4 - Kathleen Clausen
$ws.Range("A4").Value = "AB191641364"
$ws.Range("B4").Value2 = 44338.179456018515
$ws.Range("C4").Value = 5
$ws.Range("D4").Value = "unlimited"
$ws.Range("E4").Value = 4
$ws.Range("F4").Value = "Kathleen"
$ws.Range("G4").Value = "TEST 10"
$ws.Range("H4").Value = "Clausen"
$ws.Range("I4").Value = "7005 Dean Rd."
$ws.Range("J4").Value = "Indianapolis"
$ws.Range("K4").Value = "IN"
$ws.Range("L4").Value = 46220
$ws.Range("M4").Value = "United States"
$ws.Range("N4").Value = "kathleen.j.clausen@gmail.com"
$ws.Range("O4").ClearContents()
$ws.Range("P4").Value = "2.22.21.EOM1."

# --- Step 4: remove the hyperlinks (and the Hyperlink formatting they applied) ---
# The hyperlinked cells were N2/N4/N5 before the column insert in Step 1 shifted
# their formatting out to O2/O4/O5 (the "Donor Email" -> "Donor Phone" slot).
$ws.Hyperlinks.Delete()
$ws.Range("O2").ClearFormats()
$ws.Range("O4").ClearFormats()
$ws.Range("O5").ClearFormats()
$wb.Styles("Hyperlink").Delete()

# --- Step 5: clear out the old row 5 data (it becomes a blank placeholder row) ---
$ws.Range("A5:P5").ClearContents()

# --- Step 6: drop the now-superfluous trailing blank row 114 ---
$ws.Rows(114).Delete()

# --- Step 7: refresh the stored sort-range bookkeeping to match the new extent ---
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("B3:B113"))
$ws.Sort.SetRange($ws.Range("A3:P113"))
$ws.Sort.Apply()

# --- Step 8: move the active selection (cosmetic, matches the saved view state) ---
$ws.Range("E8").Select()
